$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $plo = $r - 1
    $ws.Cells.Item($r, 1).Value = "PLO$plo"
}

$ws.Cells.Item(13, 2).Value = "11" + $ws.Cells.Item(13, 2).Value()
$ws.Cells.Item(13, 3).Value = "11" + $ws.Cells.Item(13, 3).Value()

$ws.Range("F21").Select()
